# Applies the textual edits described by the diff to Instructions_jmd.docx.
#
# 1) Remove the parenthetical "(it will be brighter than the others)" from the
#    sentence about the higher-contrast grating.
# 2) Append a new sentence "Always use your index finger to press the button."
#    after the existing "...you press the right button." sentence.
# 3) Split "...with another participant (your "partner") and your goal is..."
#    into two sentences: "...(your "partner"). Your goal is..."

$d = $word.ActiveDocument

# 1) Drop the "(it will be brighter than the others)" parenthetical.
$d.Content.Find.Execute(
    "higher contrast (it will be brighter than the others). This grating",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "higher contrast. This grating", 2) | Out-Null

# 2) Add the new sentence about using the index finger.
$d.Content.Find.Execute(
    "you press the right button.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "you press the right button. Always use your index finger to press the button.", 2) | Out-Null

# 3) Split the run-on sentence into two sentences.
$d.Content.Find.Execute(
    "(your " + [char]8220 + "partner" + [char]8221 + ") and your goal is to maximize",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(your " + [char]8220 + "partner" + [char]8221 + "). Your goal is to maximize", 2) | Out-Null
